$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K" (formerly derived from Strike#). Update the
# computed K values for each data row (rows 2-14) per the regenerated
# save_data.
$newValues = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
